# "Add files via upload" - re-upload of the Internship Checklist workbook
# with a few checklist text tweaks and a couple of view/formatting touches.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Part 1 table: the two "Ask mentor to review" rows are cleared out.
$ws.Range("E10").ClearContents()
$ws.Range("E11").ClearContents()

# "Mentor visits company" row gets edited first (typo'd to "Metor visits
# company") ...
$ws.Range("E29").Value = "Metor visits company"

# ... then the "Email Supervisor & Company Details to mentor" row is
# expanded to mention Mr Kau as well.
$ws.Range("E13").Value = "Email Supervisor & Company Details to mentor and Mr Kau"

# Widen column E so the longer text fits.
$ws.Columns("E").ColumnWidth = 55.166666666666664

# Leave the selection on D29 (and scrolled back to the top of the sheet).
$ws.Range("D29").Select()
